# Update "countries & provincias Spain" COVID data sheet.
#
# The source data refreshed (new case counts for a handful of provinces),
# which shifts some rows once the table is re-sorted descending by
# "Casos totales" (column B) - exactly mirroring how the sheet was
# originally produced. The timestamp footer is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp (row 1).
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 12:55"

# Helper: look up a province/city by name in column A and overwrite its
# Casos totales / Casos activos / Recuperados / Muertes figures.
function Set-CityStats {
    param($City, $Total, $Active, $Recovered, $Deaths)

    $cell = $ws.Range("A4:A64").Find($City)
    $row = $cell.Row
    $ws.Cells.Item($row, 2).Value = $Total
    $ws.Cells.Item($row, 3).Value = $Active
    $ws.Cells.Item($row, 4).Value = $Recovered
    $ws.Cells.Item($row, 5).Value = $Deaths
}

Set-CityStats "Bizkaia/Vizcaya"     2776 1503 2010 116
Set-CityStats "Valencia/Valencia"   2508  108 2273 127
Set-CityStats "Araba/Alava"         1947 1503 1332 115
Set-CityStats "Alacant/Alicante"    1734   45 1464 225
Set-CityStats "Caceres"              991   11  896  84
Set-CityStats "Gipuzkoa/Guipuzcoa"  1017 1503  630  34
Set-CityStats "Castello/Castellon"   542    8  509  25
Set-CityStats "Badajoz"              465   40  409  16

# Re-sort the data block (rows 4-64) by Casos totales, descending - same
# as the source refresh - so rows land in their new positions.
$dataRange = $ws.Range("A4:E64")
$sortKey = $ws.Range("B4:B64")
$dataRange.Sort($sortKey, 2)
